$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header text (shared strings with volume/date) ---
$ws.Range("A8").Value = "Volume 30   Number  40"
$ws.Range("C9").Value = "Report Covering the Week  10/2/2023  Through  10/8/2023"

# --- Update crime statistics table (rows 14-30) ---
# Row 14 C14
$ws.Range("F14").Copy()
$ws.Range("C14").PasteSpecial(-4122)
$ws.Range("C14").Value = 1
# Row 14 F14
$ws.Range("F14").Value = 2
# Row 14 H14
$ws.Range("H14").Value = 100
# Row 14 I14
$ws.Range("I14").Value = 12
# Row 14 K14
$ws.Range("K14").Value = 100
# Row 14 L14
$ws.Range("L14").Value = 100
# Row 14 M14
$ws.Range("M14").Value = 9.090909090909
# Row 14 N14
$ws.Range("N14").Value = -20
# Row 15 F15
$ws.Range("F15").Value = 1
# Row 15 G15
$ws.Range("G15").Value = 2
# Row 15 H15
$ws.Range("H15").Value = -50
# Row 15 M15
$ws.Range("M15").Value = -48
# Row 15 N15
$ws.Range("N15").Value = -72.340425531914
# Row 16 C16
$ws.Range("A16").Copy()
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("C16").Value = "0"
# Row 16 E16
$ws.Range("E16").Value = -100
# Row 16 F16
$ws.Range("F16").Value = 5
# Row 16 G16
$ws.Range("G16").Value = 11
# Row 16 H16
$ws.Range("H16").Value = -54.545454545454
# Row 16 I16
$ws.Range("I16").Value = 119
# Row 16 J16
$ws.Range("J16").Value = 85
# Row 16 K16
$ws.Range("K16").Value = 40
# Row 16 L16
$ws.Range("L16").Value = 77.611940298507
# Row 16 M16
$ws.Range("M16").Value = -46.875
# Row 16 N16
$ws.Range("N16").Value = -84.625322997416
# Row 17 C17
$ws.Range("C17").Value = 5
# Row 17 D17
$ws.Range("D17").Value = 4
# Row 17 E17
$ws.Range("E17").Value = 25
# Row 17 F17
$ws.Range("F17").Value = 41
# Row 17 G17
$ws.Range("G17").Value = 23
# Row 17 H17
$ws.Range("H17").Value = 78.260869565217
# Row 17 I17
$ws.Range("I17").Value = 358
# Row 17 J17
$ws.Range("J17").Value = 305
# Row 17 K17
$ws.Range("K17").Value = 17.377049180327
# Row 17 L17
$ws.Range("L17").Value = 67.289719626168
# Row 17 M17
$ws.Range("M17").Value = 51.054852320675
# Row 17 N17
$ws.Range("N17").Value = -42.628205128205
# Row 18 D18
$ws.Range("D18").Value = 2
# Row 18 G18
$ws.Range("G18").Value = 21
# Row 18 H18
$ws.Range("H18").Value = -90.476190476190
# Row 18 J18
$ws.Range("J18").Value = 86
# Row 18 K18
$ws.Range("K18").Value = 9.302325581395
# Row 18 L18
$ws.Range("L18").Value = 95.833333333333
# Row 18 M18
$ws.Range("M18").Value = -60.669456066945
# Row 18 N18
$ws.Range("N18").Value = -93.078055964653
# Row 19 C19
$ws.Range("C19").Value = 15
# Row 19 D19
$ws.Range("D19").Value = 5
# Row 19 E19
$ws.Range("E19").Value = 200
# Row 19 F19
$ws.Range("F19").Value = 40
# Row 19 G19
$ws.Range("G19").Value = 23
# Row 19 H19
$ws.Range("H19").Value = 73.913043478260
# Row 19 I19
$ws.Range("I19").Value = 328
# Row 19 J19
$ws.Range("J19").Value = 265
# Row 19 K19
$ws.Range("K19").Value = 23.773584905660
# Row 19 L19
$ws.Range("L19").Value = 62.376237623762
# Row 19 M19
$ws.Range("M19").Value = 2.5
# Row 19 N19
$ws.Range("N19").Value = -25.791855203619
# Row 20 D20
$ws.Range("D20").Value = 3
# Row 20 E20
$ws.Range("E20").Value = 33.333333333333
# Row 20 G20
$ws.Range("G20").Value = 13
# Row 20 H20
$ws.Range("H20").Value = -7.692307692307
# Row 20 I20
$ws.Range("I20").Value = 121
# Row 20 J20
$ws.Range("J20").Value = 89
# Row 20 K20
$ws.Range("K20").Value = 35.955056179775
# Row 20 L20
$ws.Range("L20").Value = 116.071428571429
# Row 20 M20
$ws.Range("M20").Value = -21.935483870967
# Row 20 N20
$ws.Range("N20").Value = -88.206627680311
# Row 21 C21
$ws.Range("C21").Value = 25
# Row 21 D21
$ws.Range("D21").Value = 16
# Row 21 E21
$ws.Range("E21").Value = 56.25
# Row 21 F21
$ws.Range("F21").Value = 103
# Row 21 G21
$ws.Range("G21").Value = 94
# Row 21 H21
$ws.Range("H21").Value = 9.574468085106
# Row 21 I21
$ws.Range("I21").Value = 1045
# Row 21 J21
$ws.Range("J21").Value = 846
# Row 21 K21
$ws.Range("K21").Value = 23.522458628841
# Row 21 L21
$ws.Range("L21").Value = 71.592775041050
# Row 21 M21
$ws.Range("M21").Value = -13.707679603633
# Row 21 N21
$ws.Range("N21").Value = -75.618292113859
# Row 23 C23
$ws.Range("A23").Copy()
$ws.Range("C23").PasteSpecial(-4122)
$ws.Range("C23").Value = "0"
# Row 23 D23
$ws.Range("A23").Copy()
$ws.Range("D23").PasteSpecial(-4122)
$ws.Range("D23").Value = "0"
# Row 23 E23
$ws.Range("A23").Copy()
$ws.Range("E23").PasteSpecial(-4122)
$ws.Range("E23").Value = "***.*"
# Row 23 F23
$ws.Range("F23").Value = 6
# Row 23 G23
$ws.Range("G23").Value = 5
# Row 23 H23
$ws.Range("H23").Value = 20
# Row 24 C24
$ws.Range("C24").Value = 20
# Row 24 D24
$ws.Range("D24").Value = 27
# Row 24 E24
$ws.Range("E24").Value = -25.925925925925
# Row 24 F24
$ws.Range("F24").Value = 96
# Row 24 G24
$ws.Range("G24").Value = 94
# Row 24 H24
$ws.Range("H24").Value = 2.127659574468
# Row 24 I24
$ws.Range("I24").Value = 962
# Row 24 J24
$ws.Range("J24").Value = 913
# Row 24 K24
$ws.Range("K24").Value = 5.366922234392
# Row 24 L24
$ws.Range("L24").Value = 57.189542483660
# Row 24 M24
$ws.Range("M24").Value = -19.091673675357
# Row 25 C25
$ws.Range("C25").Value = 16
# Row 25 D25
$ws.Range("D25").Value = 10
# Row 25 E25
$ws.Range("E25").Value = 60
# Row 25 F25
$ws.Range("F25").Value = 65
# Row 25 G25
$ws.Range("G25").Value = 51
# Row 25 H25
$ws.Range("H25").Value = 27.450980392156
# Row 25 I25
$ws.Range("I25").Value = 564
# Row 25 J25
$ws.Range("J25").Value = 503
# Row 25 K25
$ws.Range("K25").Value = 12.127236580516
# Row 25 L25
$ws.Range("L25").Value = 32.084309133489
# Row 25 M25
$ws.Range("M25").Value = -42.448979591836
# Row 26 D26
$ws.Range("F26").Copy()
$ws.Range("D26").PasteSpecial(-4122)
$ws.Range("D26").Value = 1
# Row 26 E26
$ws.Range("H26").Copy()
$ws.Range("E26").PasteSpecial(-4122)
$ws.Range("E26").Value = -100
# Row 26 F26
$ws.Range("F26").Value = 1
# Row 26 G26
$ws.Range("G26").Value = 3
# Row 26 H26
$ws.Range("H26").Value = -66.666666666666
# Row 26 J26
$ws.Range("J26").Value = 26
# Row 26 K26
$ws.Range("K26").Value = -26.923076923076
# Row 27 F27
$ws.Range("F27").Value = 3
# Row 27 H27
$ws.Range("H27").Value = -25
# Row 27 I27
$ws.Range("I27").Value = 71
# Row 27 J27
$ws.Range("J27").Value = 57
# Row 27 K27
$ws.Range("K27").Value = 24.561403508771
# Row 27 L27
$ws.Range("L27").Value = 73.170731707317
# Row 28 C28
$ws.Range("G28").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("C28").Value = 1
# Row 28 D28
$ws.Range("A28").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("D28").Value = "0"
# Row 28 E28
$ws.Range("A28").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("E28").Value = "***.*"
# Row 28 F28
$ws.Range("G28").Copy()
$ws.Range("F28").PasteSpecial(-4122)
$ws.Range("F28").Value = 1
# Row 28 G28
$ws.Range("G28").Value = 3
# Row 28 H28
$ws.Range("H28").Value = -66.666666666666
# Row 28 I28
$ws.Range("I28").Value = 20
# Row 28 K28
$ws.Range("K28").Value = -13.043478260869
# Row 28 L28
$ws.Range("L28").Value = -13.043478260869
# Row 28 M28
$ws.Range("M28").Value = -20
# Row 28 N28
$ws.Range("N28").Value = -76.190476190476
# Row 29 C29
$ws.Range("G29").Copy()
$ws.Range("C29").PasteSpecial(-4122)
$ws.Range("C29").Value = 1
# Row 29 D29
$ws.Range("A29").Copy()
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("D29").Value = "0"
# Row 29 E29
$ws.Range("A29").Copy()
$ws.Range("E29").PasteSpecial(-4122)
$ws.Range("E29").Value = "***.*"
# Row 29 F29
$ws.Range("G29").Copy()
$ws.Range("F29").PasteSpecial(-4122)
$ws.Range("F29").Value = 1
# Row 29 G29
$ws.Range("G29").Value = 2
# Row 29 H29
$ws.Range("H29").Value = -50
# Row 29 I29
$ws.Range("I29").Value = 18
# Row 29 K29
$ws.Range("K29").Value = 5.882352941176
# Row 29 L29
$ws.Range("L29").Value = -14.285714285714
# Row 29 M29
$ws.Range("M29").Value = -21.739130434782
# Row 29 N29
$ws.Range("N29").Value = -75
# Row 30 C30
$ws.Range("A30").Copy()
$ws.Range("C30").PasteSpecial(-4122)
$ws.Range("C30").Value = "0"
